$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Data": append 6 new weekly WRESBAL observations (rows 104-109)
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

# Copy the date-cell formatting of the last existing row (A103) down over the
# new rows before filling in values, so the new date cells pick up the same
# style (number format / border / alignment) as the rest of the A column.
$wsData.Range("A103").Copy()
$wsData.Range("A104:A109").PasteSpecial(-4122)

$wsData.Range("A104").Value2 = 45189
$wsData.Range("B104").Value2 = 3231.649

$wsData.Range("A105").Value2 = 45196
$wsData.Range("B105").Value2 = 3170.324

$wsData.Range("A106").Value2 = 45203
$wsData.Range("B106").Value2 = 3145.72

$wsData.Range("A107").Value2 = 45210
$wsData.Range("B107").Value2 = 3288.945

$wsData.Range("A108").Value2 = 45217
$wsData.Range("B108").Value2 = 3353.881

$wsData.Range("A109").Value2 = 45224
$wsData.Range("B109").Value2 = 3261.886

# ---------------------------------------------------------------------------
# Sheet "SeriesInfo": refresh the series metadata pulled from FRED
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

function Set-TextValue($range, $text) {
    # Force the cell to be written as text (rather than being auto-parsed
    # into a date serial, etc.), then drop the resulting formatting so the
    # cell ends up with no explicit style - matching a plain text cell.
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.ClearFormats()
}

Set-TextValue $wsInfo.Range("B3") "2023-10-27"
Set-TextValue $wsInfo.Range("B4") "2023-10-27"
Set-TextValue $wsInfo.Range("B7") "2023-10-25"
Set-TextValue $wsInfo.Range("B14") "2023-10-26 15:35:02-05"

$wsInfo.Range("B15").Value2 = 73
